$d = $word.ActiveDocument

# --- Step 1: Append "Must support encryption (critical)." as a new list item ---
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Must support encryption (critical)."

# --- Step 2: Append "Must support compression (not critical)." as a new list item ---
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Must support compression (not critical)."

# --- Step 3: Append "Behavior" heading ---
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Behavior"
$p3.Range.ListFormat.ConvertNumbersToText()
$p3.Style = $d.Styles.Item("Heading 1")

# --- Step 4: Append the final "GonzoNet will not log..." list item with a new bullet list (numId 2) ---
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.Text = "GonzoNet will not log anything, in order to maximize speed. Instead, it will throw exceptions that clients can catch()."
$p4.Style = $d.Styles.Item("List Paragraph")

$bulletGallery = $word.ListGalleries.Item(1)
$template = $bulletGallery.ListTemplates.Item(1)
$p4.Range.ListFormat.ApplyListTemplate($template)

$appliedTemplate = $p4.Range.ListFormat.ListTemplate
$lvl0 = $appliedTemplate.ListLevels.Item(1)
$lvl0.NumberFormat = "-"

Write-Host "Done editing document."
